$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

Set-TextValue "D2" "41.274.11"
Set-TextValue "E2" "  -1.80%  "

Set-TextValue "D3" "2.175.45"
Set-TextValue "E3" "  -1.77%  "

Set-TextValue "E4" "  -0.05%  "

Set-TextValue "D5" "236.72"
Set-TextValue "E5" "  -2.03%  "

Set-TextValue "D6" "0.613"
Set-TextValue "E6" "  -2.01%  "

Set-TextValue "D7" "70.18"
Set-TextValue "E7" "  -4.34%  "

Set-TextValue "E8" "  -0.08%  "

Set-TextValue "D9" "0.581"
Set-TextValue "E9" "  -4.55%  "

Set-TextValue "D10" "40.46"
Set-TextValue "E10" "  -6.44%  "

Set-TextValue "D11" "0.0929"
Set-TextValue "E11" "  -2.84%  "

Set-TextValue "D12" "54.21"
Set-TextValue "E12" "  -5.70%  "

Set-TextValue "B13" "Polkadot"
Set-TextValue "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "6.79"
Set-TextValue "E13" "  -4.74%  "

Set-TextValue "B14" "TRON"
Set-TextValue "C14" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D14" "0.101"
Set-TextValue "E14" "  -2.12%  "

Set-TextValue "D15" "2.495.27"
Set-TextValue "E15" "  -2.02%  "

Set-TextValue "D16" "14.00"
Set-TextValue "E16" "  -1.83%  "

Set-TextValue "B17" "Polygon"
Set-TextValue "C17" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D17" "0.803"
Set-TextValue "E17" "  -4.69%  "

Set-TextValue "B18" "WrappedEther"
Set-TextValue "C18" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D18" "2.173.38"
Set-TextValue "E18" "  -1.13%  "

Set-TextValue "D19" "41.031.10"
Set-TextValue "E19" "  -2.06%  "

Set-TextValue "D20" "0.0000102"
Set-TextValue "E20" "  -7.01%  "

Set-TextValue "D21" "70.53"
Set-TextValue "E21" "  -3.26%  "

Set-TextValue "D22" "5.95"
Set-TextValue "E22" "  -3.29%  "

Set-TextValue "D23" "9.89"
Set-TextValue "E23" "  -4.95%  "

Set-TextValue "D24" "226.69"
Set-TextValue "E24" "  -1.34%  "

Set-TextValue "E25" "  -6.39%  "

Set-TextValue "E26" "  +0.01%  "

Set-TextValue "D27" "10.90"
Set-TextValue "E27" "  -5.35%  "

Set-TextValue "D28" "3.54"
Set-TextValue "E28" "  -1.98%  "

Set-TextValue "E29" "  -2.80%  "

Set-TextValue "E30" "  +0.81%  "

Set-TextValue "D31" "167.89"
Set-TextValue "E31" "  +0.25%  "

Set-TextValue "D32" "19.99"
Set-TextValue "E32" "  -2.85%  "

Set-TextValue "D33" "30.65"
Set-TextValue "E33" "  +5.41%  "

Set-TextValue "D34" "0.0770"
Set-TextValue "E34" "  -2.99%  "

Set-TextValue "D35" "5.17"
Set-TextValue "E35" "  -7.37%  "

Set-TextValue "E36" "  -2.91%  "

Set-TextValue "D37" "0.103"
Set-TextValue "E37" "  -6.30%  "

Set-TextValue "D38" "4.12"
Set-TextValue "E38" "  -3.43%  "

Set-TextValue "D39" "0.0286"
Set-TextValue "E39" "  -5.28%  "

Set-TextValue "B40" "LidoDAOToken"
Set-TextValue "C40" "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue "D40" "2.08"
Set-TextValue "E40" "  -1.63%  "

Set-TextValue "B41" "Celestia"
Set-TextValue "C41" "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D41" "11.87"
Set-TextValue "E41" "  -7.51%  "

Set-TextValue "D42" "5.43"
Set-TextValue "E42" "  -3.56%  "

Set-TextValue "D43" "59.88"
Set-TextValue "E43" "  -9.54%  "

Set-TextValue "D44" "0.192"
Set-TextValue "E44" "  -3.99%  "

Set-TextValue "D45" "0.0978"
Set-TextValue "E45" "  -2.74%  "

Set-TextValue "D46" "8.28"
Set-TextValue "E46" "  -4.79%  "

Set-TextValue "D47" "98.33"
Set-TextValue "E47" "  -5.87%  "

Set-TextValue "E48" "  -2.25%  "

Set-TextValue "E49" "  -2.28%  "

Set-TextValue "E50" "  -7.08%  "

Set-TextValue "E51" "  -2.75%  "
